$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.156.78'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -2.93%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.012.04'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.94%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '530.73'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.34%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '133.91'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.24%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.005.24'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.94%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.497'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.69%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -3.73%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -1.74%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.446'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -1.55%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000220'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.57%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.18'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.57%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.510.13'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.63%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.110'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.36%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.224.83'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -2.83%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.019.03'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.77%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.61'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.61%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '465.72'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -3.47%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.24'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.88%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.678'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -2.41%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.97'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -2.33%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '79.44'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.55%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.11'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.01%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.998'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.17%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.67'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.29%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.83'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -3.34%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.32%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.90'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +1.17%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '25.63'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.56%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +2.64%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.51'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +2.58%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '55.68'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.28'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -3.68%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.89'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.89%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '459.06'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -4.67%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.221.36'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +4.33%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0785'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.48%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0385'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -2.39%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +2.19%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.16'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.55%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '27.55'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +12.39%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.48'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -6.05%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.09%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.246'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -2.76%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.00'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -1.28%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '119.37'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -1.74%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.108'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.24%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0₃0494'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -8.73%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +7.82%  '
